$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: rename "iaest-dimension:<x>" metadata to "iaest-measure:<x>" for the
# columns that moved from dimension to measure.
$ws.Range("A2").Value = "iaest-measure:temporalidad"
$ws.Range("G2").Value = "iaest-measure:case-when-bonificacioncontrato-bonif"
$ws.Range("I2").Value = "iaest-measure:mes-nombre"
$ws.Range("L2").Value = "iaest-measure:sexo"

# Row 3: "dim" -> "medida" for the same columns
$ws.Range("A3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("L3").Value = "medida"

# Row 4: "skos:Concept" -> "xsd:int" for the same columns
$ws.Range("A4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"

# Row 5: clear the mapping-file cells for columns that are no longer dimensions
$ws.Range("A5").Value = $null
$ws.Range("G5").Value = $null
$ws.Range("I5").Value = $null
$ws.Range("L5").Value = $null
